# Auto-generated Excel COM-interop script to apply market-data refresh edits
# to the Leve profit tracker workbook (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2091.7273
$ws.Range("I40").Value = 2374.75
$ws.Range("J40").Value = 1930
$ws.Range("K40").Value = 2374.75
$ws.Range("L40").Value = 1930
$ws.Range("M40").Value = -2199.75
$ws.Range("N40").Value = -2280

$ws.Range("H45").Value = 5000
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 5000
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 15000
$ws.Range("N45").Value = -15384

$ws.Range("M64").ClearContents()
$ws.Range("H64").Value = 3999
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 3999
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 3999
$ws.Range("N64").Value = -4495

$ws.Range("M67").ClearContents()
$ws.Range("H67").Value = 3999
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 3999
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 3999
$ws.Range("N67").Value = -5715

$ws.Range("H70").Value = 7417.1665
$ws.Range("I70").Value = 10251
$ws.Range("J70").Value = 6000.25
$ws.Range("K70").Value = 30753
$ws.Range("L70").Value = 18000.75
$ws.Range("M70").Value = -30483
$ws.Range("N70").Value = -18540.75

$ws.Range("H73").Value = 7417.1665
$ws.Range("I73").Value = 10251
$ws.Range("J73").Value = 6000.25
$ws.Range("K73").Value = 30753
$ws.Range("L73").Value = 18000.75
$ws.Range("M73").Value = -29817
$ws.Range("N73").Value = -19872.75

$ws.Range("H74").Value = 3824.25
$ws.Range("I74").Value = 3766
$ws.Range("J74").Value = 3999
$ws.Range("K74").Value = 3766
$ws.Range("L74").Value = 3999
$ws.Range("M74").Value = -2830
$ws.Range("N74").Value = -5871

$ws.Range("H77").Value = 3824.25
$ws.Range("I77").Value = 3766
$ws.Range("J77").Value = 3999
$ws.Range("K77").Value = 18830
$ws.Range("L77").Value = 19995
$ws.Range("M77").Value = -14150
$ws.Range("N77").Value = -29355

$ws.Range("M103").ClearContents()
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 0

$ws.Range("H106").Value = 2475
$ws.Range("I106").Value = 2475
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 2475
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -1844

$ws.Range("H132").Value = 2905.75
$ws.Range("I132").Value = 3287.1
$ws.Range("J132").Value = 999
$ws.Range("K132").Value = 9861.299999999999
$ws.Range("L132").Value = 2997
$ws.Range("M132").Value = -7331.299999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 25000
$ws.Range("I43").Value = 25000
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 25000
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -24687

$ws.Range("H63").Value = 4363.8
$ws.Range("I63").Value = 4625.778
$ws.Range("J63").Value = 2006
$ws.Range("K63").Value = 4625.778
$ws.Range("L63").Value = 2006
$ws.Range("M63").Value = -3939.778
$ws.Range("N63").Value = -3378

$ws.Range("H66").Value = 4363.8
$ws.Range("I66").Value = 4625.778
$ws.Range("J66").Value = 2006
$ws.Range("K66").Value = 23128.89
$ws.Range("L66").Value = 10030
$ws.Range("M66").Value = -19696.89
$ws.Range("N66").Value = -16894

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 28000
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 28000
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 28000
$ws.Range("N68").Value = -29622

$ws.Range("H69").Value = 17000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 17000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 17000
$ws.Range("N69").Value = -18622

$ws.Range("H71").Value = 28000
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 28000
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 84000
$ws.Range("N71").Value = -92112

$ws.Range("H72").Value = 17000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 17000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 51000
$ws.Range("N72").Value = -59112

$ws.Range("H76").Value = 29345.5
$ws.Range("I76").Value = 26999
$ws.Range("J76").Value = 30518.75
$ws.Range("K76").Value = 26999
$ws.Range("L76").Value = 30518.75
$ws.Range("M76").Value = -26684
$ws.Range("N76").Value = -31148.75

$ws.Range("H79").Value = 29345.5
$ws.Range("I79").Value = 26999
$ws.Range("J79").Value = 30518.75
$ws.Range("K79").Value = 26999
$ws.Range("L79").Value = 30518.75
$ws.Range("M79").Value = -25907
$ws.Range("N79").Value = -32702.75

$ws.Range("H82").Value = 12851
$ws.Range("I82").Value = 5688.75
$ws.Range("J82").Value = 41500
$ws.Range("K82").Value = 5688.75
$ws.Range("L82").Value = 41500
$ws.Range("M82").Value = -5305.75
$ws.Range("N82").Value = -42266

$ws.Range("H85").Value = 12851
$ws.Range("I85").Value = 5688.75
$ws.Range("J85").Value = 41500
$ws.Range("K85").Value = 5688.75
$ws.Range("L85").Value = 41500
$ws.Range("M85").Value = -4362.75
$ws.Range("N85").Value = -44152

$ws.Range("N105").ClearContents()
$ws.Range("H105").Value = 3998
$ws.Range("I105").Value = 3998
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 3998
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -2251

$ws.Range("H107").Value = 836.5
$ws.Range("I107").Value = 803.8
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 803.8
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1116.2
$ws.Range("N107").Value = -4840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 12346
$ws.Range("I38").Value = 13519
$ws.Range("J38").Value = 10000
$ws.Range("K38").Value = 13519
$ws.Range("L38").Value = 10000
$ws.Range("M38").Value = -13142
$ws.Range("N38").Value = -10754

$ws.Range("H46").Value = 12346
$ws.Range("I46").Value = 13519
$ws.Range("J46").Value = 10000
$ws.Range("K46").Value = 13519
$ws.Range("L46").Value = 10000
$ws.Range("M46").Value = -13308
$ws.Range("N46").Value = -10422

$ws.Range("N62").ClearContents()
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0

$ws.Range("N65").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0

$ws.Range("H68").Value = 41075
$ws.Range("I68").Value = 40000
$ws.Range("J68").Value = 41433.332
$ws.Range("K68").Value = 40000
$ws.Range("L68").Value = 41433.332
$ws.Range("M68").Value = -39251
$ws.Range("N68").Value = -42931.332

$ws.Range("H71").Value = 41075
$ws.Range("I71").Value = 40000
$ws.Range("J71").Value = 41433.332
$ws.Range("K71").Value = 120000
$ws.Range("L71").Value = 124299.996
$ws.Range("M71").Value = -116256
$ws.Range("N71").Value = -131787.996

$ws.Range("H86").Value = 8666.666999999999
$ws.Range("I86").Value = 1000
$ws.Range("J86").Value = 12500
$ws.Range("K86").Value = 1000
$ws.Range("L86").Value = 12500
$ws.Range("M86").Value = 123
$ws.Range("N86").Value = -14746

$ws.Range("H89").Value = 8666.666999999999
$ws.Range("I89").Value = 1000
$ws.Range("J89").Value = 12500
$ws.Range("K89").Value = 5000
$ws.Range("L89").Value = 62500
$ws.Range("M89").Value = 616
$ws.Range("N89").Value = -73732

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 2999
$ws.Range("I63").Value = 2999
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 8997
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -8248

$ws.Range("H66").Value = 2999
$ws.Range("I66").Value = 2999
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 26991
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -23247

$ws.Range("M107").ClearContents()
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 6620000
$ws.Range("I11").Value = 8250000
$ws.Range("J11").Value = 100000
$ws.Range("K11").Value = 8250000
$ws.Range("L11").Value = 100000
$ws.Range("M11").Value = -8249861
$ws.Range("N11").Value = -100278

$ws.Range("H70").Value = 8459.200000000001
$ws.Range("I70").Value = 5071.75
$ws.Range("J70").Value = 22009
$ws.Range("K70").Value = 5071.75
$ws.Range("L70").Value = 22009
$ws.Range("M70").Value = -4801.75

$ws.Range("H73").Value = 8459.200000000001
$ws.Range("I73").Value = 5071.75
$ws.Range("J73").Value = 22009
$ws.Range("K73").Value = 5071.75
$ws.Range("L73").Value = 22009
$ws.Range("M73").Value = -4135.75

$ws.Range("H97").Value = 657.2222
$ws.Range("I97").Value = 603
$ws.Range("J97").Value = 725
$ws.Range("K97").Value = 603
$ws.Range("L97").Value = 725
$ws.Range("M97").Value = -107

$ws.Range("N106").ClearContents()
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 0

$ws.Range("H132").Value = 2987
$ws.Range("I132").Value = 2730.5
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 8191.5
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -5661.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2891.1667
$ws.Range("I22").Value = 3783
$ws.Range("J22").Value = 1999.3334
$ws.Range("K22").Value = 3783
$ws.Range("L22").Value = 1999.3334
$ws.Range("M22").Value = -3488
$ws.Range("N22").Value = -2589.3334

$ws.Range("H27").Value = 2891.1667
$ws.Range("I27").Value = 3783
$ws.Range("J27").Value = 1999.3334
$ws.Range("K27").Value = 3783
$ws.Range("L27").Value = 1999.3334
$ws.Range("M27").Value = -3676
$ws.Range("N27").Value = -2213.3334

$ws.Range("H46").Value = 4909.6
$ws.Range("I46").Value = 5933.3335
$ws.Range("J46").Value = 3374
$ws.Range("K46").Value = 5933.3335
$ws.Range("L46").Value = 3374
$ws.Range("M46").Value = -5745.3335
$ws.Range("N46").Value = -3750

$ws.Range("H55").Value = 2998.75
$ws.Range("I55").Value = 995
$ws.Range("J55").Value = 3666.6667
$ws.Range("K55").Value = 995
$ws.Range("L55").Value = 3666.6667
$ws.Range("M55").Value = -822
$ws.Range("N55").Value = -4012.6667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 22059.75
$ws.Range("I63").Value = 9995
$ws.Range("J63").Value = 26081.334
$ws.Range("K63").Value = 9995
$ws.Range("L63").Value = 26081.334
$ws.Range("M63").Value = -9371
$ws.Range("N63").Value = -27329.334

$ws.Range("H66").Value = 22059.75
$ws.Range("I66").Value = 9995
$ws.Range("J66").Value = 26081.334
$ws.Range("K66").Value = 29985
$ws.Range("L66").Value = 78244.00199999999
$ws.Range("M66").Value = -26865
$ws.Range("N66").Value = -84484.00199999999
